# Update the date values in column K (rows 10-21) from "12/27" to "12/29"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 10; $row -le 21; $row++) {
    $cell = $ws.Range("K$row")
    if ($cell.Text -eq "12/27") {
        $cell.Value2 = "12/29"
    }
}
